# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E22) is re-sorted from newest-to-oldest
# to oldest-to-newest (2212, 2301, 2302, 2303, 2304, 2305, 2306), and the
# "Valor Mora" amounts in column F for the first and last period rows
# (F16 and F22) are swapped to match the new ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order the "Periodo Mora" labels in column E (rows 16-22) so that
# periods run chronologically from oldest (2212) to newest (2306).
$ws.Range("E16").Value = "2212"
$ws.Range("E17").Value = "2301"
$ws.Range("E18").Value = "2302"
$ws.Range("E19").Value = "2303"
$ws.Range("E20").Value = "2304"
$ws.Range("E21").Value = "2305"
$ws.Range("E22").Value = "2306"

# Swap the "Valor Mora" amounts that go with the first/last periods so
# the partial-month value (29333) now travels with period 2212 and the
# full-month value (40000) travels with period 2306.
$ws.Range("F16").Value = 40000
$ws.Range("F22").Value = 29333
